# Update "想去人数" (interested-people count) figures that changed between
# data refreshes of the generated gh-pages output.
#
# Sheet "展览" (sheet1 / rId1) and sheet "全部类型" (sheet4 / rId4) both list
# the same events, so each value is updated in both sheets.

$wb = $excel.ActiveWorkbook

$exhibition = $wb.Worksheets.Item("展览")
$allTypes   = $wb.Worksheets.Item("全部类型")

# Row in "展览" -> Row in "全部类型", and the new value for column F.
$updates = @(
    @{ Row1 = 6;  Row4 = 9;  Value = 2999 },
    @{ Row1 = 8;  Row4 = 11; Value = 2016 },
    @{ Row1 = 12; Row4 = 16; Value = 954 },
    @{ Row1 = 19; Row4 = 23; Value = 7260 },
    @{ Row1 = 21; Row4 = 25; Value = 1985 },
    @{ Row1 = 24; Row4 = 29; Value = 164 },
    @{ Row1 = 26; Row4 = 31; Value = 502 },
    @{ Row1 = 33; Row4 = 37; Value = 1123 }
)

foreach ($u in $updates) {
    $exhibition.Range("F" + $u.Row1).Value = $u.Value
    $allTypes.Range("F" + $u.Row4).Value = $u.Value
}
